$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# --- Row 2: replace existing candidate "Chethan" with "Sriharshini" ---
$ws.Range("A2").Value = "Sriharshini"
$ws.Range("B2").Value = "peddisriharshini@gmail.com"
$ws.Range("E2").Value = "Hyderabad,Bangalore"

# --- Row 3: add new candidate "Chandrika" ---
$ws.Range("B3").Value = "chandrikagollashetti@gmail.com"
$ws.Range("A3").Value = "Chandrika"

# --- Domains column for both rows ---
$ws.Range("C2").Value = "Fresher,Software Engineer, Software Developer, Software Tester, UiPath"
$ws.Range("C3").Value = "Fresher,Software Engineer, Software Developer, Software Tester"

# --- YearOfExperience column (dropdown list values) for both rows ---
$ws.Range("D2").Value = "Fresher"
$ws.Range("D3").Value = "Fresher"

# --- Location for the new row ---
$ws.Range("E3").Value = "Hyderabad,Bangalore"

# Fix B2 formatting to match the rest of the row (bordered style), same as A2
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null

# Re-apply B2 value, since PasteSpecial(formats) shouldn't touch it, but make sure
$ws.Range("B2").Value = "peddisriharshini@gmail.com"

# Apply row 2's formatting to row 3 so every column keeps the matching style
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Re-apply row 3 values (PasteSpecial with formats only should not alter these, but ensure correctness)
$ws.Range("A3").Value = "Chandrika"
$ws.Range("B3").Value = "chandrikagollashetti@gmail.com"
$ws.Range("C3").Value = "Fresher,Software Engineer, Software Developer, Software Tester"
$ws.Range("D3").Value = "Fresher"
$ws.Range("E3").Value = "Hyderabad,Bangalore"

$ws.Range("E6").Select() | Out-Null
